$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Workbook window size change
$excel.Width = 13680
$excel.Height = 14620

# 2) Add the 4 new header strings (row 1) in columns CS:CV
$ws.Range("CS1").Value = "enhance_equip_beyond_num"
$ws.Range("CT1").Value = "enhance_equip_no_money"
$ws.Range("CU1").Value = "merge_equip_no_money"
$ws.Range("CV1").Value = "level_up_in_battle"

# 3) Add the 4 new Chinese strings (row 2) in columns CS:CV
$ws.Range("CS2").Value = "今天铁匠已经累啦，明天再说吧"
$ws.Range("CT2").Value = "强化装备可是个体力活，没钱俺可不干"
$ws.Range("CU2").Value = "合成装备是个时间活，时间就是金钱，朋友"
$ws.Range("CV2").Value = "在战斗中渡劫升级乃是兵家大忌"

# 4) Update sheet view: selection / top-left cell
$ws.Activate()
$ws.Range("CS1").Select()
$excel.ActiveWindow.ScrollRow = 1
$excel.ActiveWindow.ScrollColumn = 97
$ws.Range("CW9").Select()
